# Add new columns I0 (I) and IF (J) to the worksheet, filling in header and data values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - copy style from existing header cell H1 (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-67 for columns I and J
$iValues = @(4,6,8,7,7,7,8,6,7,7,6,6,10,9,9,8,9,9,9,9,10,9,9,8,9,9,9,10,9,9,9,8,9,9,8,8,9,9,8,8,8,9,9,8,8,9,9,10,8,9,8,8,9,9,8,10,9,7,8,9,8,9,9,5,4,5)
$jValues = @(6,6,8,7,7,7,8,6,7,7,6,6,10,9,9,8,9,9,9,9,10,9,9,8,9,9,10,10,9,9,9,8,9,9,9,8,9,9,8,9,9,9,9,8,9,9,9,10,9,9,8,9,9,9,8,10,9,8,8,9,8,9,9,5,4,5)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
